$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.410745620727539
$ws.Range("B1").Value = 3.507086038589478
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.730283498764038
$ws.Range("E1").Value = 3.222220420837402
